$wb = $excel.ActiveWorkbook

# ---------- Sheet "Main Code" ----------
$ws1 = $wb.Worksheets.Item("Main Code")

# Widen column D
$ws1.Columns.Item(4).ColumnWidth = 13.5546875

# New row 19 - copy formatting from an existing similar row (row 16) cell by cell,
# then overwrite with the new content so styles match exactly.
$ws1.Range("A16").Copy($ws1.Range("A19"))
$ws1.Range("A19").Value = 18

$ws1.Range("B16").Copy($ws1.Range("B19"))
$ws1.Range("B19").Value = "Correct euler angles so we do not have a jump, not good  for my att stab (yaw case)"

$ws1.Range("C16").Copy($ws1.Range("C19"))
$ws1.Range("C19").Value = "2019-07-09"

$ws1.Range("D13").Copy($ws1.Range("D19"))
$ws1.Range("D19").Value = "BUG"

$ws1.Range("F16").Copy($ws1.Range("F19"))
$ws1.Range("F19").Value = "ToDo"

$ws1.Rows.Item(19).RowHeight = 28.8

# ---------- Sheet "Visualization Code" ----------
$ws2 = $wb.Worksheets.Item("Visualization Code")
$ws2.Range("G7").Value = "Make propellers move "

# ---------- back to "Main Code": row 20 ----------
$ws1.Range("A16").Copy($ws1.Range("A20"))
$ws1.Range("A20").Value = 19

$ws1.Range("B16").Copy($ws1.Range("B20"))
$ws1.Range("B20").Value = "Attitude stabilization with PID"

$ws1.Range("C16").Copy($ws1.Range("C20"))
$ws1.Range("C20").Value = "2019-07-09"

$ws1.Range("D4").Copy($ws1.Range("D20"))
$ws1.Range("D20").Value = "DEV"

$ws1.Range("E15").Copy($ws1.Range("E20"))
$ws1.Range("E20").Value = "2019-07-10"

$ws1.Range("F16").Copy($ws1.Range("F20"))
$ws1.Range("F20").Value = "Done"

# D15 / D16 / D18 get "IMPROVEMENT" (new shared string, added last)
$ws1.Range("D15").Value = "IMPROVEMENT"
$ws1.Range("D16").Value = "IMPROVEMENT"
$ws1.Range("D18").Value = "IMPROVEMENT"

# Sheet view: scroll back to top-left, select D16
$ws1.Application.Goto($ws1.Range("A1"))
$ws1.Range("D16").Select()
$ws2.Range("G7").Select()
